$wb = $excel.ActiveWorkbook

$booksSheet = $wb.Worksheets.Item("NYT Books")
$boxSheet   = $wb.Worksheets.Item("Box Office")

# --- Add the two new worksheets at the end, in order ---
$weekly = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $boxSheet)
$weekly.Name = "Weekly Data"

$corr = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weekly)
$corr.Name = "Correlation"

# =========================================================
# Weekly Data sheet
# =========================================================

# Header row (string entry order matters for sharedStrings index order)
$weekly.Range("B1").Value = "Week End"
$weekly.Range("A1").Value = "Week Start"
$weekly.Range("C1").Value = "Rank"
$weekly.Range("D1").Value = "Box Office Gross"

# Column formatting: column B (Week End) uses the existing date format,
# column A/B per-cell date style, column D width for "Box Office Gross"
$weekly.Columns.Item(2).NumberFormat = "yyyy\-mm\-dd;@"
$weekly.Columns.Item(2).ColumnWidth = 13.71
$weekly.Columns.Item(4).ColumnWidth = 14.0

# Week End values (column B) -- one per week, Saturday dates
$weekEnds = @(40251,40258,40265,40272,40279,40286,40293,40300,40307,40314,40321,40328,40335)
for ($i = 0; $i -lt $weekEnds.Count; $i++) {
    $row = $i + 2
    $weekly.Range("B$row").Value = $weekEnds[$i]
    $weekly.Range("B$row").NumberFormat = "yyyy\-mm\-dd;@"
}

# Week Start formulas (column A) = Week End - 6, shared across A3:A14
$weekly.Range("A2").Formula = "=B2-6"
$weekly.Range("A2").NumberFormat = "yyyy\-mm\-dd;@"
$weekly.Range("A3:A14").Formula = "=B3-6"
$weekly.Range("A3:A14").NumberFormat = "yyyy\-mm\-dd;@"

# Rank values (column C), copied from "NYT Books" column B
$ranks = @(10,8,8,6,6,5,5,7,6,6,10,9,10)
for ($i = 0; $i -lt $ranks.Count; $i++) {
    $row = $i + 2
    $weekly.Range("C$row").Value = $ranks[$i]
}

# Box Office Gross formulas (column D), SUM of the matching Box Office rows
$weekly.Range("D4").Formula  = "=SUM('Box Office'!B2:B4)"
$weekly.Range("D5").Formula  = "=SUM('Box Office'!B5:B11)"
$weekly.Range("D6").Formula  = "=SUM('Box Office'!B12:B18)"
$weekly.Range("D7").Formula  = "=SUM('Box Office'!B19:B25)"
$weekly.Range("D8").Formula  = "=SUM('Box Office'!B26:B32)"
$weekly.Range("D9").Formula  = "=SUM('Box Office'!B33:B39)"
$weekly.Range("D10").Formula = "=SUM('Box Office'!B40:B46)"
$weekly.Range("D11").Formula = "=SUM('Box Office'!B47:B53)"
$weekly.Range("D12").Formula = "=SUM('Box Office'!B54:B60)"
$weekly.Range("D13").Formula = "=SUM('Box Office'!B61:B67)"
$weekly.Range("D14").Formula = "=SUM('Box Office'!B68:B74)"

# =========================================================
# Correlation sheet
# =========================================================

$corr.Range("A1").Value = "Rank"
$corr.Range("B1").Value = "Box Office Gross"
$corr.Columns.Item(2).ColumnWidth = 14.0

$corrRanks  = @(8,6,6,5,5,7,6,6,10,9,10)
$grossVals = @(43732319,48403597,41268522,24846628,20094861,13827823,8840117,6633829,3343222,1764135,1128090)
for ($i = 0; $i -lt $corrRanks.Count; $i++) {
    $row = $i + 2
    $corr.Range("A$row").Value = $corrRanks[$i]
    $corr.Range("B$row").Value = $grossVals[$i]
}

$corr.Range("D2").Value = "Correlation:"
$corr.Range("D3").Formula = "=CORREL(A2:A12,B2:B12)"

# =========================================================
# View / selection state to match the target workbook
# =========================================================

$booksSheet.Range("B1:D14").Select()
$boxSheet.Select()
$boxSheet.Range("C23").Select()
$weekly.Range("C4:D14").Select()
$corr.Range("D4").Select()
$corr.Select()
